$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing hours entry for week 14-20/12/2015 (row 11, 8th report)
$ws.Range("B11").Formula = "=3"

# Move the active selection to B12 as left by the editor
$ws.Range("B12").Select()
